$wb = $excel.ActiveWorkbook

# Update "想去人数" (F column) counts on both the "展览" sheet and the
# "全部类型" sheet, which carry duplicate listings for the same events.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F3").Value = 1430
    $ws.Range("F4").Value = 96
    $ws.Range("F6").Value = 14
}
